$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table (rows 2-172) holds one price observation per row, sorted by
# date. A new weekly observation is inserted as row 136, which pushes the
# previously-existing rows 136-172 down to 137-173 (dimension grows from
# A1:R172 to A1:R173). Insert a whole row at position 136 so formatting of
# subsequent rows shifts down intact.
$ws.Rows.Item(136).Insert()

# Populate the newly inserted row 136 with the new observation.
$ws.Cells.Item(136, 1).Value = 10
$ws.Cells.Item(136, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(136, 3).Value = "La Araucanía"
$ws.Cells.Item(136, 4).Value = 44736
$ws.Cells.Item(136, 5).Value = 9
$ws.Cells.Item(136, 6).Value = 100114007
$ws.Cells.Item(136, 7).Value = "Jengibre"
$ws.Cells.Item(136, 8).Value = "Sin especificar"
$ws.Cells.Item(136, 9).Value = "Primera"
$ws.Cells.Item(136, 10).Value = 40
$ws.Cells.Item(136, 11).Value = 20000
$ws.Cells.Item(136, 12).Value = 20000
$ws.Cells.Item(136, 13).Value = 20000
$ws.Cells.Item(136, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(136, 15).Value = "Perú"
$ws.Cells.Item(136, 16).Value = 1538
$ws.Cells.Item(136, 17).Value = 13
$ws.Cells.Item(136, 18).Value = "Hortaliza"
